$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 511, shifting existing rows 511:608 down to 512:609.
$ws.Rows.Item(511).Insert()

# Populate the newly inserted row 511 with its data.
$ws.Range("A511").Value = 4
$ws.Range("B511").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C511").Value = "Los Lagos"
$ws.Range("D511").Value = 45211
$ws.Range("E511").Value = 10
$ws.Range("F511").Value = 100112023
$ws.Range("G511").Value = "Brócoli"
$ws.Range("H511").Value = "Sin especificar"
$ws.Range("I511").Value = "Primera"
$ws.Range("J511").Value = 500
$ws.Range("K511").Value = 1500
$ws.Range("L511").Value = 1500
$ws.Range("M511").Value = 1500
$ws.Range("N511").Value = "`$/unidad"
$ws.Range("O511").Value = "Región Metropolitana"
$ws.Range("P511").Value = 1500
$ws.Range("Q511").Value = 1
$ws.Range("R511").Value = "Hortaliza"
